$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-15 (existing rows 6..13 of the "extr" series) get renamed/shifted,
# and two new rows (16,17) are appended at the bottom.
# Final state for rows 8..17, columns A..E:

$data = @(
    @(8,  6,  "line7", 14, 11, $false),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# New rows (16 and 17) need the same bold/bordered/centered style that
# column A already uses for every other data row; copy it from A15.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "done"
